{"js": "// Helper: search for an exact literal string and replace the single match.\nasync function replaceExact(oldText, newText) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText.substring(0, 60));\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1. Update \"Ativa\u00e7\u00e3o\" date\nawait replaceExact(\"Ativa\u00e7\u00e3o: 01/01/1996\", \"Ativa\u00e7\u00e3o: 01/01/2022\");\n\n// 2. Replace the \"Objetivos\" paragraph text\nawait replaceExact(\n  \"Fornecer aos estudantes uma vis\u00e3o abrangente e interdisciplinar dos materiais compostos por fases caracterizadas por distintos tipos de materiais (metais, cer\u00e2micas e pol\u00edmeros) para obter propriedades \u00fanicas. Apresentar os fundamentos te\u00f3ricos da mec\u00e2nica de estruturas refor\u00e7adas com fibras, tecidos e part\u00edculas. Apresentar os diferentes tipos de materiais comp\u00f3sitos, inclusive sobre os nanocomp\u00f3sitos e comp\u00f3sitos funcionais, que representam o avan\u00e7o mais recente na \u00e1rea de Ci\u00eancia e Engenharia de Materiais.\",\n  \"Fornecer aos estudantes uma vis\u00e3o abrangente e interdisciplinar sobre materiais comp\u00f3sitos, al\u00e9m de mostrar as especificidades de cada matriz, sendo ela met\u00e1lica, cer\u00e2mica ou polim\u00e9rica. Ademais, deseja-se apresentar os fundamentos te\u00f3ricos da mec\u00e2nica de estruturas refor\u00e7adas e a partir de atividades pr\u00e1ticas demostrar m\u00e9todos de caracteriza\u00e7\u00e3o de materiais comp\u00f3sitos e como prepara-los.\"\n);\n\n// 3. Add two new professors after \"519033 - Carlos Yujiro Shigue\"\nconst paras = context.document.body.paragraphs;\nparas.load(\"text\");\nawait context.sync();\nlet profPara = null;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text === \"519033 - Carlos Yujiro Shigue\") {\n    profPara = paras.items[i];\n    break;\n  }\n}\nprofPara.insertBreak(Word.BreakType.line, Word.InsertLocation.end);\nawait context.sync();\nprofPara.insertText(\"1033242 - F\u00e1bio Herbst Florenzano\", Word.InsertLocation.end);\nawait context.sync();\nprofPara.insertBreak(Word.BreakType.line, Word.InsertLocation.end);\nawait context.sync();\nprofPara.insertText(\"1922320 - Sebastiao Ribeiro\", Word.InsertLocation.end);\nawait context.sync();\n\n// 4. Replace \"Programa resumido\" paragraph text\nconst r4 = context.document.body.search(\n  \"Materiais comp\u00f3sitos: tipos, propriedades, processamento e aplica\u00e7\u00f5es. Nanocomp\u00f3sitos e comp\u00f3sitos funcionais.\",\n  { matchCase: true }\n);\nawait context.sync();\nr4.items[0].insertText(\n  \"1.Introdu\u00e7\u00e2o 2. Conceitos b\u00e1sicos sobre materiais comp\u00f3sitos, suas matrizes e seus processo de fabrica\u00e7\u00e3o 3. Tipos de refor\u00e7os 4. Comp\u00f3sitos nanoestruturados, naturais e h\u00edbridos 5. Mec\u00e2nica da estrutura refor\u00e7ada 6. Atividade pr\u00e1tica\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// 5. Replace \"Programa\" paragraph text\nconst r5 = context.document.body.search(\n  \"Conte\u00fado te\u00f3rico:1. Conceitos b\u00e1sicos sobre materiais comp\u00f3sitos: comp\u00f3sitos de matriz met\u00e1lica (CMM), comp\u00f3sitos de matriz cer\u00e2micos (CMC) e comp\u00f3sitos de matriz polim\u00e9rica (CMP) e nanocomp\u00f3sitos.2. Fibras, tecidos e refor\u00e7os particulados.3. Mec\u00e2nica de estruturas refor\u00e7adas.4. Comp\u00f3sitos de matriz met\u00e1lica: caracter\u00edsticas e processos de fabrica\u00e7\u00e3o.5. Comp\u00f3sitos de matriz cer\u00e2mica: caracter\u00edsticas e processos de fabrica\u00e7\u00e3o.6. Comp\u00f3sitos de matriz polim\u00e9rica: matrizes termopl\u00e1sticas e termorr\u00edgidas, caracter\u00edsticas f\u00edsicas e qu\u00edmicas e processos de fabrica\u00e7\u00e3o.7. Comp\u00f3sitos nanoestruturados.8.Comp\u00f3sitos funcionais.Conte\u00fado pr\u00e1tico:1. Caracteriza\u00e7\u00e3o e an\u00e1lise de comp\u00f3sitos de matriz met\u00e1lica.2. Prepara\u00e7\u00e3o e caracteriza\u00e7\u00e3o de comp\u00f3sito de matriz polim\u00e9rica.3. Visita a empresa produtora de comp\u00f3sitos.\",\n  { matchCase: true }\n);\nawait context.sync();\nr5.items[0].insertText(\n  \"1. Conceitos b\u00e1sicos sobre materiais comp\u00f3sitos: comp\u00f3sitos de matriz met\u00e1lica (CMM), comp\u00f3sitos de matriz cer\u00e2micos (CMC) e comp\u00f3sitos de matriz polim\u00e9rica (CMP) e nanocomp\u00f3sitos. 2. Tipos de Refor\u00e7os: Refor\u00e7os particulados, fibras curtas, fibras longas, mantas, tecidos e preformas. 3. Conceitos de Interface4. Comp\u00f3sitos de matriz met\u00e1lica: caracter\u00edsticas e processos de fabrica\u00e7\u00e3o. 5. Comp\u00f3sitos de matriz cer\u00e2mica: caracter\u00edsticas e processos de fabrica\u00e7\u00e3o. 6. Comp\u00f3sitos de matriz polim\u00e9rica: matrizes termopl\u00e1sticas e termorr\u00edgidas, caracter\u00edsticas f\u00edsicas e qu\u00edmicas e processos de fabrica\u00e7\u00e3o. 7. Comp\u00f3sitos nanoestruturados. 8. Comp\u00f3sitos Naturais. 9. Comp\u00f3sitos H\u00edbridos 10. Mec\u00e2nica de estruturas refor\u00e7adas. Conte\u00fado pr\u00e1tico: 1. Caracteriza\u00e7\u00e3o e an\u00e1lise de comp\u00f3sitos de matriz met\u00e1lica. 2. Prepara\u00e7\u00e3o e caracteriza\u00e7\u00e3o de comp\u00f3sitos de matriz polim\u00e9rica.(Sugest\u00e3o: Considerar substituir essa parte pr\u00e1tica pela realiza\u00e7\u00e3o do PBL descrito no item 3) 3. Visita a empresa produtora de comp\u00f3sitos e aulas especiais e/ou palestras com professores/pesquisadores convidados\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// 6. Replace \"M\u00e9todo:\" text\nconst r6 = context.document.body.search(\n  \"A avalia\u00e7\u00e3o ser\u00e1 feita por meio de provas escritas.\",\n  { matchCase: true }\n);\nawait context.sync();\nr6.items[0].insertText(\n  \"De acordo com a atual ementa da disciplina prop\u00f5e-se o uso de uma nova metodologia de ensino com o intuito de abordar o conte\u00fado de forma mais pr\u00e1tica e contextualizada para que o aluno consiga relacionar os conhecimentos te\u00f3ricos vistos em sala de aula com as outras disciplinas do curso. Assim, avalia\u00e7\u00e3o do aluno ser\u00e1 feita atrav\u00e9s de uma prova escrita e por uma apresenta\u00e7\u00e3o final com base nas atividades pr\u00e1ticas desenvolvidas.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// 7. Replace \"Crit\u00e9rio:\" text\nconst r7 = context.document.body.search(\n  \"A Nota final (NF) ser\u00e1 calculada da seguinte maneira:NF = (P1 + 2*P2)/3\",\n  { matchCase: true }\n);\nawait context.sync();\nr7.items[0].insertText(\n  \"A nota final ser\u00e1 calculada como descrita a seguir: NF= (0,4*Avalia\u00e7\u00e3o escrita + 0,6 *Apresenta\u00e7\u00e3o final)\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// 8. Replace \"Norma de recupera\u00e7\u00e3o:\" text\nconst r8 = context.document.body.search(\n  \"A recupera\u00e7\u00e3o ser\u00e1 feita por meio de uma prova escrita (PR) e a m\u00e9dia de recupera\u00e7\u00e3o (MR) calculada pela f\u00f3rmula: MR = (NF + PR)/2\",\n  { matchCase: true }\n);\nawait context.sync();\nr8.items[0].insertText(\n  \"Devido a cunho pr\u00e1tico da disciplina n\u00e3o haver\u00e1 recupera\u00e7\u00e3o.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// 9. Replace the Bibliografia paragraph text\nconst r9 = context.document.body.search(\n  \"1. MALLICK, P.K. Composites Engineering Handbook. New York: Marcel Dekker, 1997.2. MATTHEWS, F.L. & RAWLINGS, R.D. Composite Materials: Engineering and Science. London: Chapman & Hall, 1994.3. OBRAZTSOV, I.F. Mechanics of Composites. Moscow: MIR Publishers, 1982.4. JONES R. Mechanics of Composite Materials. New York: McGraw-Hill, 1975.5. UPADHYAYA, G.S. Sintered Metal-Ceramic Composites. Elsevier, 1984.6. HARPER, C. A. Handbook of Plastics, Elastomers and Composites. New York: McGraw-Hill, 1992.7. GOLDSTEIN, A.N. Handbook of Nanophase Materials. CRC Press, 1997.8. DRESSELHAUS, M.S. Graphite Fibers and Filaments. New York: Springer-Verlag, 1988.\",\n  { matchCase: true }\n);\nawait context.sync();\nr9.items[0].insertText(\n  \"1. REZENDE, M. C.; COSTA, M. L.; BOTELHO, E. C. Comp\u00f3sitos estruturais: tecnologia e pr\u00e1tica. S\u00e3o Paulo: Artliber, 2011. 396p. 2 MALLICK, P.K. Composites Engineering Handbook. New York: Marcel Dekker, 1997. 3. MATTHEWS, F.L. & RAWLINGS, R.D. Composite Materials: Engineering and Science. London: Chapman & Hall, 1994. 4. OBRAZTSOV, I.F. Mechanics of Composites. Moscow: MIR Publishers, 1982. 5. JONES R. Mechanics of Composite Materials. New York: McGraw-Hill, 1975. 6. UPADHYAYA, G.S. Sintered Metal-Ceramic Composites. Elsevier, 1984. 7. HARPER, C. A. Handbook of Plastics, Elastomers and Composites. New York: McGraw-Hill, 1992. 8. GOLDSTEIN, A.N. Handbook of Nanophase Materials. CRC Press, 1997. 9. DRESSELHAUS, M.S. Graphite Fibers and Filaments. New York: Springer-Verlag, 1988.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-ExactText($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.MatchWildcards = $false\n    $find.MatchCase = $true\n    $result = $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, [Microsoft.Office.Interop.Word.WdReplace]::wdReplaceOne)\n    if (-not $result) {\n        throw \"Replace failed for text starting: $($oldText.Substring(0, [Math]::Min(60, $oldText.Length)))\"\n    }\n}\n\n# 1. Update \"Ativa\u00e7\u00e3o\" date\nReplace-ExactText \"Ativa\u00e7\u00e3o: 01/01/1996\" \"Ativa\u00e7\u00e3o: 01/01/2022\"\n\n# 2. Replace the \"Objetivos\" paragraph text\nReplace-ExactText \"Fornecer aos estudantes uma vis\u00e3o abrangente e interdisciplinar dos materiais compostos por fases caracterizadas por distintos tipos de materiais (metais, cer\u00e2micas e pol\u00edmeros) para obter propriedades \u00fanicas. Apresentar os fundamentos te\u00f3ricos da mec\u00e2nica de estruturas refor\u00e7adas com fibras, tecidos e part\u00edculas. Apresentar os diferentes tipos de materiais comp\u00f3sitos, inclusive sobre os nanocomp\u00f3sitos e comp\u00f3sitos funcionais, que representam o avan\u00e7o mais recente na \u00e1rea de Ci\u00eancia e Engenharia de Materiais.\" `\n    \"Fornecer aos estudantes uma vis\u00e3o abrangente e interdisciplinar sobre materiais comp\u00f3sitos, al\u00e9m de mostrar as especificidades de cada matriz, sendo ela met\u00e1lica, cer\u00e2mica ou polim\u00e9rica. Ademais, deseja-se apresentar os fundamentos te\u00f3ricos da mec\u00e2nica de estruturas refor\u00e7adas e a partir de atividades pr\u00e1ticas demostrar m\u00e9todos de caracteriza\u00e7\u00e3o de materiais comp\u00f3sitos e como prepara-los.\"\n\n# 3. Add two new professors after \"519033 - Carlos Yujiro Shigue\"\n$findProf = $d.Content.Find\n$findProf.ClearFormatting()\n$findProf.Text = \"519033 - Carlos Yujiro Shigue\"\n$findProf.MatchWildcards = $false\n$findProf.MatchCase = $true\n$foundProf = $findProf.Execute()\nif (-not $foundProf) {\n    throw \"Could not find professor anchor paragraph\"\n}\n$profRange = $findProf.Parent\n$profRange.Collapse([Microsoft.Office.Interop.Word.WdCollapseDirection]::wdCollapseEnd)\n$profRange.InsertAfter(\"`v1033242 - F\u00e1bio Herbst Florenzano`v1922320 - Sebastiao Ribeiro\")\n\n# 4. Replace \"Programa resumido\" paragraph text\nReplace-ExactText \"Materiais comp\u00f3sitos: tipos, propriedades, processamento e aplica\u00e7\u00f5es. Nanocomp\u00f3sitos e comp\u00f3sitos funcionais.\" `\n    \"1.Introdu\u00e7\u00e2o 2. Conceitos b\u00e1sicos sobre materiais comp\u00f3sitos, suas matrizes e seus processo de fabrica\u00e7\u00e3o 3. Tipos de refor\u00e7os 4. Comp\u00f3sitos nanoestruturados, naturais e h\u00edbridos 5. Mec\u00e2nica da estrutura refor\u00e7ada 6. Atividade pr\u00e1tica\"\n\n# 5. Replace \"Programa\" paragraph text\nReplace-ExactText \"Conte\u00fado te\u00f3rico:1. Conceitos b\u00e1sicos sobre materiais comp\u00f3sitos: comp\u00f3sitos de matriz met\u00e1lica (CMM), comp\u00f3sitos de matriz cer\u00e2micos (CMC) e comp\u00f3sitos de matriz polim\u00e9rica (CMP) e nanocomp\u00f3sitos.2. Fibras, tecidos e refor\u00e7os particulados.3. Mec\u00e2nica de estruturas refor\u00e7adas.4. Comp\u00f3sitos de matriz met\u00e1lica: caracter\u00edsticas e processos de fabrica\u00e7\u00e3o.5. Comp\u00f3sitos de matriz cer\u00e2mica: caracter\u00edsticas e processos de fabrica\u00e7\u00e3o.6. Comp\u00f3sitos de matriz polim\u00e9rica: matrizes termopl\u00e1sticas e termorr\u00edgidas, caracter\u00edsticas f\u00edsicas e qu\u00edmicas e processos de fabrica\u00e7\u00e3o.7. Comp\u00f3sitos nanoestruturados.8.Comp\u00f3sitos funcionais.Conte\u00fado pr\u00e1tico:1. Caracteriza\u00e7\u00e3o e an\u00e1lise de comp\u00f3sitos de matriz met\u00e1lica.2. Prepara\u00e7\u00e3o e caracteriza\u00e7\u00e3o de comp\u00f3sito de matriz polim\u00e9rica.3. Visita a empresa produtora de comp\u00f3sitos.\" `\n    \"1. Conceitos b\u00e1sicos sobre materiais comp\u00f3sitos: comp\u00f3sitos de matriz met\u00e1lica (CMM), comp\u00f3sitos de matriz cer\u00e2micos (CMC) e comp\u00f3sitos de matriz polim\u00e9rica (CMP) e nanocomp\u00f3sitos. 2. Tipos de Refor\u00e7os: Refor\u00e7os particulados, fibras curtas, fibras longas, mantas, tecidos e preformas. 3. Conceitos de Interface4. Comp\u00f3sitos de matriz met\u00e1lica: caracter\u00edsticas e processos de fabrica\u00e7\u00e3o. 5. Comp\u00f3sitos de matriz cer\u00e2mica: caracter\u00edsticas e processos de fabrica\u00e7\u00e3o. 6. Comp\u00f3sitos de matriz polim\u00e9rica: matrizes termopl\u00e1sticas e termorr\u00edgidas, caracter\u00edsticas f\u00edsicas e qu\u00edmicas e processos de fabrica\u00e7\u00e3o. 7. Comp\u00f3sitos nanoestruturados. 8. Comp\u00f3sitos Naturais. 9. Comp\u00f3sitos H\u00edbridos 10. Mec\u00e2nica de estruturas refor\u00e7adas. Conte\u00fado pr\u00e1tico: 1. Caracteriza\u00e7\u00e3o e an\u00e1lise de comp\u00f3sitos de matriz met\u00e1lica. 2. Prepara\u00e7\u00e3o e caracteriza\u00e7\u00e3o de comp\u00f3sitos de matriz polim\u00e9rica.(Sugest\u00e3o: Considerar substituir essa parte pr\u00e1tica pela realiza\u00e7\u00e3o do PBL descrito no item 3) 3. Visita a empresa produtora de comp\u00f3sitos e aulas especiais e/ou palestras com professores/pesquisadores convidados\"\n\n# 6. Replace \"M\u00e9todo:\" text\nReplace-ExactText \"A avalia\u00e7\u00e3o ser\u00e1 feita por meio de provas escritas.\" `\n    \"De acordo com a atual ementa da disciplina prop\u00f5e-se o uso de uma nova metodologia de ensino com o intuito de abordar o conte\u00fado de forma mais pr\u00e1tica e contextualizada para que o aluno consiga relacionar os conhecimentos te\u00f3ricos vistos em sala de aula com as outras disciplinas do curso. Assim, avalia\u00e7\u00e3o do aluno ser\u00e1 feita atrav\u00e9s de uma prova escrita e por uma apresenta\u00e7\u00e3o final com base nas atividades pr\u00e1ticas desenvolvidas.\"\n\n# 7. Replace \"Crit\u00e9rio:\" text\nReplace-ExactText \"A Nota final (NF) ser\u00e1 calculada da seguinte maneira:NF = (P1 + 2*P2)/3\" `\n    \"A nota final ser\u00e1 calculada como descrita a seguir: NF= (0,4*Avalia\u00e7\u00e3o escrita + 0,6 *Apresenta\u00e7\u00e3o final)\"\n\n# 8. Replace \"Norma de recupera\u00e7\u00e3o:\" text\nReplace-ExactText \"A recupera\u00e7\u00e3o ser\u00e1 feita por meio de uma prova escrita (PR) e a m\u00e9dia de recupera\u00e7\u00e3o (MR) calculada pela f\u00f3rmula: MR = (NF + PR)/2\" `\n    \"Devido a cunho pr\u00e1tico da disciplina n\u00e3o haver\u00e1 recupera\u00e7\u00e3o.\"\n\n# 9. Replace the Bibliografia paragraph text\nReplace-ExactText \"1. MALLICK, P.K. Composites Engineering Handbook. New York: Marcel Dekker, 1997.2. MATTHEWS, F.L. & RAWLINGS, R.D. Composite Materials: Engineering and Science. London: Chapman & Hall, 1994.3. OBRAZTSOV, I.F. Mechanics of Composites. Moscow: MIR Publishers, 1982.4. JONES R. Mechanics of Composite Materials. New York: McGraw-Hill, 1975.5. UPADHYAYA, G.S. Sintered Metal-Ceramic Composites. Elsevier, 1984.6. HARPER, C. A. Handbook of Plastics, Elastomers and Composites. New York: McGraw-Hill, 1992.7. GOLDSTEIN, A.N. Handbook of Nanophase Materials. CRC Press, 1997.8. DRESSELHAUS, M.S. Graphite Fibers and Filaments. New York: Springer-Verlag, 1988.\" `\n    \"1. REZENDE, M. C.; COSTA, M. L.; BOTELHO, E. C. Comp\u00f3sitos estruturais: tecnologia e pr\u00e1tica. S\u00e3o Paulo: Artliber, 2011. 396p. 2 MALLICK, P.K. Composites Engineering Handbook. New York: Marcel Dekker, 1997. 3. MATTHEWS, F.L. & RAWLINGS, R.D. Composite Materials: Engineering and Science. London: Chapman & Hall, 1994. 4. OBRAZTSOV, I.F. Mechanics of Composites. Moscow: MIR Publishers, 1982. 5. JONES R. Mechanics of Composite Materials. New York: McGraw-Hill, 1975. 6. UPADHYAYA, G.S. Sintered Metal-Ceramic Composites. Elsevier, 1984. 7. HARPER, C. A. Handbook of Plastics, Elastomers and Composites. New York: McGraw-Hill, 1992. 8. GOLDSTEIN, A.N. Handbook of Nanophase Materials. CRC Press, 1997. 9. DRESSELHAUS, M.S. Graphite Fibers and Filaments. New York: Springer-Verlag, 1988.\"\n"}
